$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (H1) to new header cells I1, J1
$ws.Range("H1").Copy($ws.Range("I1:J1"))

# Header values for new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0)
$colI = @(6, 2, 9, 7, 5, 4, 4, 4)
# Data values for column J (IF)
$colJ = @(9, 5, 9, 8, 7, 6, 7, 5)

for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
